# Trade #16 (row 17) closed out. Update the Summary, Strategy Status,
# All Trades, and MarketMaking sheets to reflect the closed trade.

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.01   # Current Capital
$summary.Range("B4").Value = 0.01      # Total P&L $
$summary.Range("B5").Value = 0.01      # Total P&L %
$summary.Range("B6").Value = 16        # Total Trades
$summary.Range("B7").Value = 5         # Winning Trades
$summary.Range("B9").Value = 31.25     # Win Rate %

# ---- Strategy Status sheet (MarketMaking row, row 4) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.01     # Capital
$status.Range("D4").Value = 16         # Trades
$status.Range("E4").Value = 0.01       # P&L $
$status.Range("F4").Value = 0.01       # P&L %
$status.Range("G4").Value = 31.25      # Win Rate %

# ---- All Trades sheet (Trade #16, row 17) ----
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G17").Value = 0.97
$allTrades.Range("H17").Value = "CLOSED"
$allTrades.Range("I17").Value = 83.0189
$allTrades.Range("J17").Value = 0.44
$allTrades.Range("K17").Value = 100.01
$allTrades.Range("P17").Value = "early_exit"
$allTrades.Range("Q17").Value = 5.02

# ---- MarketMaking sheet (Trade #16, row 17) ----
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G17").Value = 0.97
$mm.Range("H17").Value = "CLOSED"
$mm.Range("I17").Value = 83.0189
$mm.Range("J17").Value = 0.44
$mm.Range("K17").Value = 100.01
$mm.Range("P17").Value = "early_exit"
$mm.Range("Q17").Value = 5.02
